$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Mark these watchlist test cases as not running ("N"), leaving C6 untouched
$ws.Range("C5").Value = "N"
$ws.Range("C7").Value = "N"
$ws.Range("C8").Value = "N"
$ws.Range("C9").Value = "N"
$ws.Range("C10").Value = "N"
$ws.Range("C11").Value = "N"

# Update the active selection to C6, the test case that is still set to run
$ws.Activate()
$ws.Range("C6").Select()
